# FRED WALCL data refresh: append latest weekly observations and update
# the SeriesInfo metadata to match the newly pulled series.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

# --- Append two new weekly observations to the Data sheet ---------------
# Copy the formatting of the last existing data row down onto the new
# rows before filling in the actual date / value pairs.
$wsData.Range("A97").Copy()
$wsData.Range("A98").PasteSpecial(-4122)
$wsData.Range("A99").PasteSpecial(-4122)

$wsData.Range("A98").Value = 45147
$wsData.Range("B98").Value = 8208.241

$wsData.Range("A99").Value = 45154
$wsData.Range("B99").Value = 8145.727

# --- Refresh SeriesInfo metadata to match the newer FRED pull -----------
# Leading apostrophes force these date-looking values to stay as plain
# text (matching the inlineStr cells already used on this sheet) instead
# of being auto-converted to date serial numbers.
$wsInfo.Range("B3").Value = "'2023-08-22"
$wsInfo.Range("B4").Value = "'2023-08-22"
$wsInfo.Range("B7").Value = "'2023-08-16"
$wsInfo.Range("B14").Value = "'2023-08-17 15:33:36-05"
